$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5: round data values (columns B..AH) to 2 decimal places ("custom accuracy").
# Column A (timestamp) and columns already at 2dp (P5, Y5) are left as-is.
$row5 = @{
    "B5"  = 19.7
    "C5"  = 14.46
    "D5"  = 1.21
    "E5"  = 42.82
    "F5"  = 34.85
    "G5"  = 15.48
    "H5"  = 61.35
    "I5"  = 23.85
    "J5"  = 10.56
    "K5"  = 15.59
    "L5"  = 17.18
    "M5"  = 18.1
    "N5"  = 4.95
    "O5"  = 15.41
    "P5"  = 21.91
    "Q5"  = 13.05
    "R5"  = 0.82
    "S5"  = 0.82
    "T5"  = 227.32
    "U5"  = 43.14
    "V5"  = 14.23
    "W5"  = 28.93
    "X5"  = 15.14
    "Y5"  = 2.36
    "Z5"  = 29.62
    "AA5" = 12.57
    "AB5" = 11.16
    "AC5" = 13.14
    "AD5" = 18.04
    "AE5" = 0.55
    "AF5" = 55.87
    "AG5" = 7.97
    "AH5" = 17.79
}

foreach ($addr in $row5.Keys) {
    $ws.Range($addr).Value = $row5[$addr]
}

# Remove the last data row (row 6) entirely, shrinking the used range to A1:AH5.
$ws.Rows(6).Delete()
